# Update non-optimizer ascended data
#
# The "Optimizer Disabled - Ascended" sheet has its column D (mint cost)
# inputs lowered for rows 3-13. Columns E, L, M, Q, R on that sheet are
# formulas (=C+D, =E-J, =L-L(prev), etc.) so they recompute automatically
# once D is written. The two embedded charts on this workbook plot the
# Q/R columns of this same sheet, so their cached series follow along once
# Excel recalculates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optimizer Disabled - Ascended")

$ws.Range("D3").Value  = 84992
$ws.Range("D4").Value  = 368915
$ws.Range("D5").Value  = 729510
$ws.Range("D6").Value  = 1137605
$ws.Range("D7").Value  = 1593200
$ws.Range("D8").Value  = 2096295
$ws.Range("D9").Value  = 2646890
$ws.Range("D10").Value = 3244985
$ws.Range("D11").Value = 3890580
$ws.Range("D12").Value = 4583675
$ws.Range("D13").Value = 5304178

# Recalculate so every dependent formula (E/L/M/Q/R columns + chart
# caches) is refreshed before the workbook is saved.
$excel.CalculateFull()

# The active sheet/selection moved: "Optimizer Disabled - Ascended" becomes
# the selected tab (it previously was "Optimizer | 1000 run - Descend"),
# and its selection moved from D25 to D15.
$ws.Activate()
$ws.Range("D15").Select()
